# Automatic update of files.
# Column C ("Förändrad") holds the last-changed date for each cleared-felling
# notice. Bump it from 2023-11-03 (serial 45233) to 2023-11-13 (serial 45243)
# for every data row (rows 2-18) on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C18").Value = 45243
